# Update column G ("K") values in the active worksheet to reflect a
# regeneration of save data that now uses K (strikeouts) instead of the
# previous Strike# derived figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 2
    4  = 7
    5  = 2
    6  = 5
    7  = 2
    8  = 0
    9  = 2
    10 = 5
    11 = 4
    12 = 4
    13 = 2
    14 = 2
    15 = 1
    16 = 3
    17 = 4
    18 = 8
    19 = 1
    20 = 0
    21 = 3
    22 = 1
    23 = 1
    24 = 4
    25 = 5
    26 = 5
    27 = 3
    28 = 2
    29 = 4
    30 = 2
    31 = 3
    32 = 3
    33 = 3
    34 = 2
    35 = 1
    36 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
